{"js": "// Map of original equation text -> replacement equation text, in the\n// same row-major order as the cells in the (single) table.\nconst replacements = [\n  [\"24-12=12\", \"90-65=25\"],\n  [\"39+23=62\", \"31+6=37\"],\n  [\"37-28=9\", \"17+44=61\"],\n  [\"5+59=64\", \"61+33=94\"],\n  [\"87-26=61\", \"41-4=37\"],\n  [\"74-3=71\", \"83-71=12\"],\n  [\"13+81=94\", \"82+0=82\"],\n  [\"18-6=12\", \"27+46=73\"],\n  [\"34-5=29\", \"85-36=49\"],\n  [\"8+43=51\", \"68-44=24\"],\n  [\"35+57=92\", \"27-7=20\"],\n  [\"79-55=24\", \"35+59=94\"],\n  [\"75-29=46\", \"43-12=31\"],\n  [\"59-15=44\", \"25+51=76\"],\n  [\"27+60=87\", \"88-74=14\"],\n  [\"34+43=77\", \"91+0=91\"],\n  [\"58+35=93\", \"41+13=54\"],\n  [\"42+15=57\", \"82-59=23\"],\n  [\"86-16=70\", \"92-87=5\"],\n  [\"64-44=20\", \"89-72=17\"],\n  [\"38+40=78\", \"75+12=87\"],\n  [\"75-0=75\", \"39-4=35\"],\n  [\"59+5=64\", \"84-49=35\"],\n  [\"20+41=61\", \"98-38=60\"],\n  [\"7+73=80\", \"32-13=19\"],\n  [\"64-63=1\", \"62-31=31\"],\n  [\"79-38=41\", \"20-15=5\"],\n  [\"92-67=25\", \"35-1=34\"],\n  [\"65+22=87\", \"13-7=6\"],\n  [\"69-44=25\", \"71-38=33\"],\n  [\"50-5=45\", \"72-59=13\"],\n  [\"16-11=5\", \"1+0=1\"],\n  [\"70-61=9\", \"48-19=29\"],\n  [\"45-10=35\", \"43+31=74\"],\n  [\"66-34=32\", \"86-46=40\"],\n  [\"52-0=52\", \"25+46=71\"],\n  [\"13+37=50\", \"46-4=42\"],\n  [\"28+52=80\", \"76-10=66\"],\n  [\"58+3=61\", \"93-87=6\"],\n  [\"14-9=5\", \"68-39=29\"],\n  [\"39+21=60\", \"25+18=43\"],\n  [\"13+30=43\", \"83-7=76\"],\n  [\"99-47=52\", \"89-76=13\"],\n  [\"11+50=61\", \"72-48=24\"],\n  [\"84-22=62\", \"0+51=51\"],\n  [\"79-29=50\", \"22+19=41\"],\n  [\"53+14=67\", \"40-36=4\"],\n  [\"71-65=6\", \"51-29=22\"],\n  [\"74+20=94\", \"52+11=63\"],\n  [\"69+4=73\", \"11+16=27\"],\n  [\"13+73=86\", \"26+27=53\"],\n  [\"87-44=43\", \"62+26=88\"],\n  [\"95-24=71\", \"56+26=82\"],\n  [\"65-30=35\", \"58-1=57\"],\n  [\"58-8=50\", \"30+21=51\"],\n  [\"16+80=96\", \"16+53=69\"],\n  [\"81-31=50\", \"32+31=63\"],\n  [\"50-6=44\", \"73-30=43\"],\n  [\"31-14=17\", \"22+32=54\"],\n  [\"53+22=75\", \"70+12=82\"],\n  [\"43-16=27\", \"22+24=46\"],\n  [\"38+0=38\", \"18+51=69\"],\n  [\"48-43=5\", \"91-26=65\"],\n  [\"9+16=25\", \"80+3=83\"],\n  [\"49-36=13\", \"63+5=68\"],\n  [\"67-41=26\", \"24+32=56\"],\n  [\"81-26=55\", \"23+72=95\"],\n  [\"64+3=67\", \"87-47=40\"],\n  [\"63+6=69\", \"53+33=86\"],\n  [\"36+35=71\", \"47+8=55\"],\n  [\"7+17=24\", \"96-90=6\"],\n  [\"56+21=77\", \"55-37=18\"],\n  [\"9+67=76\", \"57-12=45\"],\n  [\"75+16=91\", \"85-61=24\"],\n  [\"94-73=21\", \"11+84=95\"],\n  [\"14-3=11\", \"80-1=79\"],\n  [\"13+23=36\", \"6+1=7\"],\n  [\"65+17=82\", \"51-27=24\"],\n  [\"90-72=18\", \"83-59=24\"],\n  [\"63+26=89\", \"75-47=28\"],\n  [\"37+32=69\", \"15+84=99\"],\n  [\"28+15=43\", \"65-57=8\"],\n  [\"77-13=64\", \"41+19=60\"],\n  [\"11+59=70\", \"48-19=29\"],\n  [\"47-4=43\", \"62-25=37\"],\n  [\"75-47=28\", \"13+86=99\"],\n  [\"77-51=26\", \"78-35=43\"],\n  [\"76+14=90\", \"38-11=27\"],\n  [\"41-35=6\", \"94-78=16\"],\n  [\"45+30=75\", \"98-22=76\"],\n  [\"21+66=87\", \"84-11=73\"],\n  [\"21+58=79\", \"1+41=42\"],\n  [\"82-48=34\", \"99-73=26\"],\n  [\"11+60=71\", \"31-22=9\"],\n  [\"76+3=79\", \"70-20=50\"],\n  [\"66-21=45\", \"18+54=72\"],\n  [\"33+36=69\", \"87+2=89\"],\n  [\"66+24=90\", \"44-26=18\"],\n  [\"45+9=54\", \"84-38=46\"],\n  [\"57-39=18\", \"19+39=58\"],\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = table.values[0].length;\n\n// Walk the cells in row-major order (matches `replacements` order) and\n// swap in the new equation text, cell by cell, preserving each cell's\n// existing paragraph/run formatting.\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const [oldText, newText] = replacements[idx];\n\n    // Defensive check: only touch the cell if it still holds the\n    // expected original value (avoids a later replacement accidentally\n    // re-matching text that an earlier replacement just wrote).\n    const currentText = table.values[r][c].trim();\n    if (currentText === oldText) {\n      const cell = table.getCell(r, c);\n      const range = cell.body.getRange();\n      range.insertText(newText, Word.InsertLocation.replace);\n    }\n    idx++;\n  }\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Old -> new equation text, in row-major cell order (matches the table\n# layout: 20 rows x 5 columns = 100 cells).\n$pairs = @(\n    @('24-12=12', '90-65=25'),\n    @('39+23=62', '31+6=37'),\n    @('37-28=9', '17+44=61'),\n    @('5+59=64', '61+33=94'),\n    @('87-26=61', '41-4=37'),\n    @('74-3=71', '83-71=12'),\n    @('13+81=94', '82+0=82'),\n    @('18-6=12', '27+46=73'),\n    @('34-5=29', '85-36=49'),\n    @('8+43=51', '68-44=24'),\n    @('35+57=92', '27-7=20'),\n    @('79-55=24', '35+59=94'),\n    @('75-29=46', '43-12=31'),\n    @('59-15=44', '25+51=76'),\n    @('27+60=87', '88-74=14'),\n    @('34+43=77', '91+0=91'),\n    @('58+35=93', '41+13=54'),\n    @('42+15=57', '82-59=23'),\n    @('86-16=70', '92-87=5'),\n    @('64-44=20', '89-72=17'),\n    @('38+40=78', '75+12=87'),\n    @('75-0=75', '39-4=35'),\n    @('59+5=64', '84-49=35'),\n    @('20+41=61', '98-38=60'),\n    @('7+73=80', '32-13=19'),\n    @('64-63=1', '62-31=31'),\n    @('79-38=41', '20-15=5'),\n    @('92-67=25', '35-1=34'),\n    @('65+22=87', '13-7=6'),\n    @('69-44=25', '71-38=33'),\n    @('50-5=45', '72-59=13'),\n    @('16-11=5', '1+0=1'),\n    @('70-61=9', '48-19=29'),\n    @('45-10=35', '43+31=74'),\n    @('66-34=32', '86-46=40'),\n    @('52-0=52', '25+46=71'),\n    @('13+37=50', '46-4=42'),\n    @('28+52=80', '76-10=66'),\n    @('58+3=61', '93-87=6'),\n    @('14-9=5', '68-39=29'),\n    @('39+21=60', '25+18=43'),\n    @('13+30=43', '83-7=76'),\n    @('99-47=52', '89-76=13'),\n    @('11+50=61', '72-48=24'),\n    @('84-22=62', '0+51=51'),\n    @('79-29=50', '22+19=41'),\n    @('53+14=67', '40-36=4'),\n    @('71-65=6', '51-29=22'),\n    @('74+20=94', '52+11=63'),\n    @('69+4=73', '11+16=27'),\n    @('13+73=86', '26+27=53'),\n    @('87-44=43', '62+26=88'),\n    @('95-24=71', '56+26=82'),\n    @('65-30=35', '58-1=57'),\n    @('58-8=50', '30+21=51'),\n    @('16+80=96', '16+53=69'),\n    @('81-31=50', '32+31=63'),\n    @('50-6=44', '73-30=43'),\n    @('31-14=17', '22+32=54'),\n    @('53+22=75', '70+12=82'),\n    @('43-16=27', '22+24=46'),\n    @('38+0=38', '18+51=69'),\n    @('48-43=5', '91-26=65'),\n    @('9+16=25', '80+3=83'),\n    @('49-36=13', '63+5=68'),\n    @('67-41=26', '24+32=56'),\n    @('81-26=55', '23+72=95'),\n    @('64+3=67', '87-47=40'),\n    @('63+6=69', '53+33=86'),\n    @('36+35=71', '47+8=55'),\n    @('7+17=24', '96-90=6'),\n    @('56+21=77', '55-37=18'),\n    @('9+67=76', '57-12=45'),\n    @('75+16=91', '85-61=24'),\n    @('94-73=21', '11+84=95'),\n    @('14-3=11', '80-1=79'),\n    @('13+23=36', '6+1=7'),\n    @('65+17=82', '51-27=24'),\n    @('90-72=18', '83-59=24'),\n    @('63+26=89', '75-47=28'),\n    @('37+32=69', '15+84=99'),\n    @('28+15=43', '65-57=8'),\n    @('77-13=64', '41+19=60'),\n    @('11+59=70', '48-19=29'),\n    @('47-4=43', '62-25=37'),\n    @('75-47=28', '13+86=99'),\n    @('77-51=26', '78-35=43'),\n    @('76+14=90', '38-11=27'),\n    @('41-35=6', '94-78=16'),\n    @('45+30=75', '98-22=76'),\n    @('21+66=87', '84-11=73'),\n    @('21+58=79', '1+41=42'),\n    @('82-48=34', '99-73=26'),\n    @('11+60=71', '31-22=9'),\n    @('76+3=79', '70-20=50'),\n    @('66-21=45', '18+54=72'),\n    @('33+36=69', '87+2=89'),\n    @('66+24=90', '44-26=18'),\n    @('45+9=54', '84-38=46'),\n    @('57-39=18', '19+39=58')\n)\n\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $oldText = $pairs[$idx][0]\n        $newText = $pairs[$idx][1]\n        $cell = $tbl.Cell($r, $c)\n        $range = $cell.Range\n        # Cell range text includes trailing cell-mark characters; trim them off.\n        $currentText = $range.Text.TrimEnd([char]13, [char]7)\n\n        # Defensive check: only replace the cell if it still holds the\n        # expected original value (avoids a later replacement accidentally\n        # re-matching text that an earlier replacement just wrote).\n        if ($currentText -eq $oldText) {\n            $cell.Range.Text = $newText\n        }\n        $idx++\n    }\n}\n"}
